$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 208
$ws.Range("I2").Value = 209.2
$ws.Range("J2").Value = 202
$ws.Range("K2").Value = 209.2
$ws.Range("L2").Value = 202
$ws.Range("M2").Value = -96.19999999999999
$ws.Range("N2").Value = -428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9995
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 644.6923
$ws.Range("I41").Value = 524.625
$ws.Range("K41").Value = 524.625
$ws.Range("M41").Value = -84.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 26499.75
$ws.Range("J43").Value = 2000
$ws.Range("L43").Value = 2000
$ws.Range("N43").Value = -2138

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 37332.332
$ws.Range("J116").Value = 5999
$ws.Range("L116").Value = 5999
$ws.Range("N116").Value = -12883

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2570.6667
$ws.Range("I132").Value = 2570.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7712.000100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5182.000100000001
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 36347.793
$ws.Range("I137").Value = 53500.57
$ws.Range("J137").Value = 8639.462
$ws.Range("K137").Value = 160501.71
$ws.Range("L137").Value = 25918.386
$ws.Range("M137").Value = -157951.71
$ws.Range("N137").Value = -31018.386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1999.5
$ws.Range("I45").Value = 1999.25
$ws.Range("K45").Value = 1999.25
$ws.Range("M45").Value = -1622.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3582.2163
$ws.Range("I61").Value = 2081.4736
$ws.Range("J61").Value = 5166.3335
$ws.Range("K61").Value = 2081.4736
$ws.Range("L61").Value = 5166.3335
$ws.Range("M61").Value = -1869.4736
$ws.Range("N61").Value = -5590.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 8429.647000000001
$ws.Range("I110").Value = 9332.154
$ws.Range("K110").Value = 9332.154
$ws.Range("M110").Value = -7287.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3582.2163
$ws.Range("I136").Value = 2081.4736
$ws.Range("J136").Value = 5166.3335
$ws.Range("K136").Value = 6244.4208
$ws.Range("L136").Value = 15499.0005
$ws.Range("M136").Value = -3694.4208
$ws.Range("N136").Value = -20599.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 20999
$ws.Range("J15").Value = 20999
$ws.Range("L15").Value = 20999
$ws.Range("N15").Value = -21453

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2484.25
$ws.Range("I86").Value = 2320.4
$ws.Range("K86").Value = 2320.4
$ws.Range("M86").Value = -1197.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2484.25
$ws.Range("I89").Value = 2320.4
$ws.Range("K89").Value = 11602
$ws.Range("M89").Value = -5986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2563.625
$ws.Range("I16").Value = 2227.5
$ws.Range("K16").Value = 2227.5
$ws.Range("M16").Value = -1940.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 7274.3
$ws.Range("I41").Value = 2985.5
$ws.Range("K41").Value = 2985.5
$ws.Range("M41").Value = -2557.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2018.75
$ws.Range("I58").Value = 2275.1667
$ws.Range("J58").Value = 1249.5
$ws.Range("K58").Value = 2275.1667
$ws.Range("L58").Value = 1249.5
$ws.Range("M58").Value = -2072.1667
$ws.Range("N58").Value = -1655.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1288.2609
$ws.Range("I94").Value = 848.7
$ws.Range("J94").Value = 1626.3846
$ws.Range("K94").Value = 848.7
$ws.Range("L94").Value = 1626.3846
$ws.Range("M94").Value = -397.7
$ws.Range("N94").Value = -2528.3846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 447746
$ws.Range("I99").Value = 926903.2
$ws.Range("J99").Value = 42305.31
$ws.Range("K99").Value = 926903.2
$ws.Range("L99").Value = 42305.31
$ws.Range("M99").Value = -925405.2
$ws.Range("N99").Value = -45301.31

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2563.625
$ws.Range("I113").Value = 2227.5
$ws.Range("K113").Value = 2227.5
$ws.Range("M113").Value = -57.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 4993
$ws.Range("I122").Value = 4993
$ws.Range("K122").Value = 14979
$ws.Range("M122").Value = -12529

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 447746
$ws.Range("I126").Value = 926903.2
$ws.Range("J126").Value = 42305.31
$ws.Range("K126").Value = 2780709.6
$ws.Range("L126").Value = 126915.93
$ws.Range("M126").Value = -2778239.6
$ws.Range("N126").Value = -131855.93

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6314.8
$ws.Range("I134").Value = 8155.913
$ws.Range("K134").Value = 24467.739
$ws.Range("M134").Value = -21932.739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2018.75
$ws.Range("I136").Value = 2275.1667
$ws.Range("J136").Value = 1249.5
$ws.Range("K136").Value = 6825.500100000001
$ws.Range("L136").Value = 3748.5
$ws.Range("M136").Value = -4275.500100000001
$ws.Range("N136").Value = -8848.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 69.40000000000001
$ws.Range("I14").Value = 69.40000000000001
$ws.Range("K14").Value = 208.2
$ws.Range("M14").Value = -35.20000000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 4000
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 901.5
$ws.Range("I113").Value = 899.3333
$ws.Range("K113").Value = 2697.9999
$ws.Range("M113").Value = -527.9998999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10059.571
$ws.Range("I43").Value = 8899.75
$ws.Range("J43").Value = 11606
$ws.Range("K43").Value = 8899.75
$ws.Range("L43").Value = 11606
$ws.Range("M43").Value = -8748.75
$ws.Range("N43").Value = -11908

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7823.067
$ws.Range("I70").Value = 7753.96
$ws.Range("K70").Value = 7753.96
$ws.Range("M70").Value = -7483.96

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7823.067
$ws.Range("I73").Value = 7753.96
$ws.Range("K73").Value = 7753.96
$ws.Range("M73").Value = -6817.96

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5916.6665
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1966.1154
$ws.Range("I61").Value = 1999.2632
$ws.Range("J61").Value = 1876.1428
$ws.Range("K61").Value = 1999.2632
$ws.Range("L61").Value = 1876.1428
$ws.Range("M61").Value = -1797.2632
$ws.Range("N61").Value = -2280.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 8551.75
$ws.Range("I68").Value = 11360.647
$ws.Range("K68").Value = 11360.647
$ws.Range("M68").Value = -10611.647

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 8551.75
$ws.Range("I71").Value = 11360.647
$ws.Range("K71").Value = 56803.235
$ws.Range("M71").Value = -53059.235

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1966.1154
$ws.Range("I113").Value = 1999.2632
$ws.Range("J113").Value = 1876.1428
$ws.Range("K113").Value = 1999.2632
$ws.Range("L113").Value = 1876.1428
$ws.Range("M113").Value = 170.7367999999999
$ws.Range("N113").Value = -6216.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1465.5
$ws.Range("J100").Value = 1403
$ws.Range("N100").Value = -3888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1329
$ws.Range("I132").Value = 1689
$ws.Range("K132").Value = 5067
$ws.Range("M132").Value = -2537
